$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Checking out order #2: record quantities/total and mark it as processed
$ws.Range("B2").Value = "u"
$ws.Range("C2").Value = "p"
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 3217.5

# Inventory/order tracking advances to the next order number
$ws.Range("A3").Value = 3
